$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the changed Price/Volume cells keep their original text (string) data type
# rather than being auto-converted to numbers/percentages by Excel.
$textCells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "D48", "E48", "D49", "E49", "E50", "D51", "E51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "305.91"
$ws.Range("E2").Value = "0.83%"
$ws.Range("D3").Value = "37.09"
$ws.Range("E3").Value = "6.05%"
$ws.Range("D4").Value = "5.050"
$ws.Range("E4").Value = "-2.22%"
$ws.Range("D5").Value = "0.07814"
$ws.Range("E5").Value = "-0.44%"
$ws.Range("D6").Value = "2.215"
$ws.Range("E6").Value = "-4.91%"
$ws.Range("D7").Value = "8.001"
$ws.Range("E7").Value = "-0.68%"
$ws.Range("D8").Value = "4.012"
$ws.Range("E8").Value = "0.91%"
$ws.Range("D9").Value = "0.9279"
$ws.Range("D10").Value = "0.09826"
$ws.Range("E10").Value = "-2.45%"
$ws.Range("D11").Value = "0.1886"
$ws.Range("E11").Value = "3.42%"
$ws.Range("D12").Value = "0.08670"
$ws.Range("E12").Value = "1.63%"
$ws.Range("D13").Value = "0.03692"
$ws.Range("E13").Value = "8.37%"
$ws.Range("D14").Value = "0.09971"
$ws.Range("E14").Value = "0.61%"
$ws.Range("D15").Value = "0.001470"
$ws.Range("E15").Value = "-0.47%"
$ws.Range("D16").Value = "0.005681"
$ws.Range("E16").Value = "-2.29%"
$ws.Range("D17").Value = "3.462"
$ws.Range("E17").Value = "-0.28%"
$ws.Range("D18").Value = "2.333"
$ws.Range("E18").Value = "10.97%"
$ws.Range("E19").Value = "-0.46%"
$ws.Range("D20").Value = "0.1315"
$ws.Range("E20").Value = "-0.76%"
$ws.Range("D21").Value = "4.745"
$ws.Range("E21").Value = "4.81%"
$ws.Range("D22").Value = "0.2199"
$ws.Range("E22").Value = "-0.93%"
$ws.Range("D23").Value = "0.04579"
$ws.Range("E23").Value = "-1.17%"
$ws.Range("D24").Value = "0.001251"
$ws.Range("E24").Value = "2.93%"
$ws.Range("D25").Value = "0.004497"
$ws.Range("E25").Value = "0.85%"
$ws.Range("D26").Value = "0.0001398"
$ws.Range("E26").Value = "7.75%"
$ws.Range("D27").Value = "0.0002715"
$ws.Range("E27").Value = "-19.94%"
$ws.Range("D39").Value = "0.01870"
$ws.Range("E39").Value = "6.46%"
$ws.Range("D40").Value = "0.04787"
$ws.Range("E40").Value = "0.98%"
$ws.Range("D41").Value = "0.007971"
$ws.Range("E41").Value = "2.56%"
$ws.Range("D42").Value = "0.1410"
$ws.Range("E42").Value = "-0.38%"
$ws.Range("D43").Value = "0.007552"
$ws.Range("E43").Value = "-14.27%"
$ws.Range("D44").Value = "0.002101"
$ws.Range("E44").Value = "-8.02%"
$ws.Range("D45").Value = "0.01012"
$ws.Range("E45").Value = "10.38%"
$ws.Range("D46").Value = "0.00006387"
$ws.Range("E46").Value = "5.55%"
$ws.Range("D47").Value = "0.00000000749"
$ws.Range("E47").Value = "0.07%"
$ws.Range("D48").Value = "0.0005785"
$ws.Range("E48").Value = "-0.26%"
$ws.Range("D49").Value = "30.42"
$ws.Range("E49").Value = "682.95%"
$ws.Range("E50").Value = "-36.02%"
$ws.Range("D51").Value = "0.00002097"
$ws.Range("E51").Value = "0.07%"
